$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 105, shifting existing rows 105-116 down to 106-117
$ws.Rows.Item(105).Insert()

# Populate the newly inserted row 105 with the new weekly record
$ws.Cells.Item(105, 1).Value = 7
$ws.Cells.Item(105, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(105, 3).Value = "Ñuble"
$ws.Cells.Item(105, 4).Value = 45021
$ws.Cells.Item(105, 5).Value = 16
$ws.Cells.Item(105, 6).Value = "Fruta"
$ws.Cells.Item(105, 7).Value = 100108
$ws.Cells.Item(105, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(105, 9).Value = 100108002
$ws.Cells.Item(105, 10).Value = "Mango"
$ws.Cells.Item(105, 11).Value = "Sin especificar"
$ws.Cells.Item(105, 12).Value = "Primera"
$ws.Cells.Item(105, 13).Value = 50
$ws.Cells.Item(105, 14).Value = 7000
$ws.Cells.Item(105, 15).Value = 7000
$ws.Cells.Item(105, 16).Value = 7000
$ws.Cells.Item(105, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(105, 18).Value = "Perú"
$ws.Cells.Item(105, 19).Value = 1750
$ws.Cells.Item(105, 20).Value = 4
